$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - F-column (想去人数) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4949
$ws1.Range("F4").Value = 342
$ws1.Range("F7").Value = 47
$ws1.Range("F9").Value = 133
$ws1.Range("F10").Value = 316
$ws1.Range("F11").Value = 253
$ws1.Range("F12").Value = 2972
$ws1.Range("F13").Value = 154
$ws1.Range("F14").Value = 1564

# Sheet "全部类型" (sheet4) - F-column (想去人数) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4949
$ws4.Range("F4").Value = 342
$ws4.Range("F8").Value = 47
$ws4.Range("F10").Value = 133
$ws4.Range("F11").Value = 316
$ws4.Range("F12").Value = 253
$ws4.Range("F13").Value = 2972
$ws4.Range("F14").Value = 154
$ws4.Range("F15").Value = 1564
